{"js": "const replacements = [\n  [\"865\u00f73=288, 1\", \"215\u00f75=43, 0\"],\n  [\"929\u00f74=232, 1\", \"632\u00f75=126, 2\"],\n  [\"946\u00f76=157, 4\", \"889\u00f77=127, 0\"],\n  [\"330\u00f78=41, 2\", \"376\u00f78=47, 0\"],\n  [\"535\u00f76=89, 1\", \"120\u00f75=24, 0\"],\n  [\"683\u00f79=75, 8\", \"899\u00f73=299, 2\"],\n  [\"945\u00f78=118, 1\", \"440\u00f73=146, 2\"],\n  [\"122\u00f73=40, 2\", \"844\u00f78=105, 4\"],\n  [\"811\u00f78=101, 3\", \"389\u00f74=97, 1\"],\n  [\"691\u00f75=138, 1\", \"361\u00f76=60, 1\"],\n  [\"112\u00f73=37, 1\", \"255\u00f76=42, 3\"],\n  [\"114\u00f76=19, 0\", \"620\u00f73=206, 2\"],\n  [\"959\u00f74=239, 3\", \"497\u00f79=55, 2\"],\n  [\"472\u00f79=52, 4\", \"628\u00f77=89, 5\"],\n  [\"458\u00f77=65, 3\", \"236\u00f72=118, 0\"],\n  [\"589\u00f79=65, 4\", \"865\u00f74=216, 1\"],\n  [\"845\u00f73=281, 2\", \"823\u00f73=274, 1\"],\n  [\"603\u00f76=100, 3\", \"166\u00f76=27, 4\"],\n  [\"410\u00f79=45, 5\", \"525\u00f76=87, 3\"],\n  [\"974\u00f74=243, 2\", \"676\u00f73=225, 1\"],\n  [\"808\u00f79=89, 7\", \"807\u00f76=134, 3\"],\n  [\"697\u00f74=174, 1\", \"997\u00f76=166, 1\"],\n  [\"125\u00f78=15, 5\", \"804\u00f74=201, 0\"],\n  [\"138\u00f73=46, 0\", \"262\u00f76=43, 4\"],\n  [\"339\u00f75=67, 4\", \"211\u00f77=30, 1\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"865\u00f73=288, 1\", \"215\u00f75=43, 0\"),\n    @(\"929\u00f74=232, 1\", \"632\u00f75=126, 2\"),\n    @(\"946\u00f76=157, 4\", \"889\u00f77=127, 0\"),\n    @(\"330\u00f78=41, 2\", \"376\u00f78=47, 0\"),\n    @(\"535\u00f76=89, 1\", \"120\u00f75=24, 0\"),\n    @(\"683\u00f79=75, 8\", \"899\u00f73=299, 2\"),\n    @(\"945\u00f78=118, 1\", \"440\u00f73=146, 2\"),\n    @(\"122\u00f73=40, 2\", \"844\u00f78=105, 4\"),\n    @(\"811\u00f78=101, 3\", \"389\u00f74=97, 1\"),\n    @(\"691\u00f75=138, 1\", \"361\u00f76=60, 1\"),\n    @(\"112\u00f73=37, 1\", \"255\u00f76=42, 3\"),\n    @(\"114\u00f76=19, 0\", \"620\u00f73=206, 2\"),\n    @(\"959\u00f74=239, 3\", \"497\u00f79=55, 2\"),\n    @(\"472\u00f79=52, 4\", \"628\u00f77=89, 5\"),\n    @(\"458\u00f77=65, 3\", \"236\u00f72=118, 0\"),\n    @(\"589\u00f79=65, 4\", \"865\u00f74=216, 1\"),\n    @(\"845\u00f73=281, 2\", \"823\u00f73=274, 1\"),\n    @(\"603\u00f76=100, 3\", \"166\u00f76=27, 4\"),\n    @(\"410\u00f79=45, 5\", \"525\u00f76=87, 3\"),\n    @(\"974\u00f74=243, 2\", \"676\u00f73=225, 1\"),\n    @(\"808\u00f79=89, 7\", \"807\u00f76=134, 3\"),\n    @(\"697\u00f74=174, 1\", \"997\u00f76=166, 1\"),\n    @(\"125\u00f78=15, 5\", \"804\u00f74=201, 0\"),\n    @(\"138\u00f73=46, 0\", \"262\u00f76=43, 4\"),\n    @(\"339\u00f75=67, 4\", \"211\u00f77=30, 1\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $ok = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $ok) {\n        Write-Output \"FAILED to replace: $oldText\"\n    }\n}\nWrite-Output \"done\"\n"}
